# Reorder "Recorded By" (column G) names so that "System" (exact case)
# is moved to immediately follow the first token when the first token is
# the lowercase "system", or moved to the very front otherwise.
#
# Examples:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"             -> "System, backup@backdoor.com"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "system, backup@backdoor.com, System"     -> "system, System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf("System") -lt 0) { continue }

    $parts = $val -split ', '
    if ($parts.Count -lt 2) { continue }

    $sysIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $sysIndex = $i
            break
        }
    }
    if ($sysIndex -lt 0) { continue }

    # Determine target index: right after first element if first element
    # is the lowercase literal "system", otherwise the very front (index 0).
    $targetIndex = 0
    if ($parts[0].Equals("system")) {
        $targetIndex = 1
    }

    if ($sysIndex -eq $targetIndex) { continue }

    # Build new order: remove "System" from its current spot, insert at target
    $remaining = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $sysIndex) {
            $remaining += $parts[$i]
        }
    }

    $newPartsList = @()
    for ($i = 0; $i -lt $remaining.Count; $i++) {
        if ($i -eq $targetIndex) { $newPartsList += "System" }
        $newPartsList += $remaining[$i]
    }
    if ($targetIndex -ge $remaining.Count) { $newPartsList += "System" }

    $newVal = [string]::Join(", ", $newPartsList)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
